$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Simple cell value updates (price / volume changes)
Set-TextValue "D2" "25.746.46"
$ws.Range("E2").Value = "  -0.04%  "
Set-TextValue "D3" "1.756.21"
$ws.Range("E3").Value = "  -1.79%  "
$ws.Range("E4").Value = "  -0.70%  "
Set-TextValue "D5" "236.74"
$ws.Range("E5").Value = "  -1.62%  "
Set-TextValue "D6" "0.9998"
$ws.Range("E6").Value = "  -0.55%  "
Set-TextValue "D7" "0.5054"
$ws.Range("E7").Value = "  +0.10%  "
Set-TextValue "D8" "41.15"
$ws.Range("E8").Value = "  -4.35%  "
Set-TextValue "D9" "0.2647"
$ws.Range("E9").Value = "  +6.25%  "
Set-TextValue "D10" "0.06189"
$ws.Range("E10").Value = "  +1.59%  "
Set-TextValue "D11" "1.747.43"
$ws.Range("E11").Value = "  -2.40%  "
Set-TextValue "D12" "0.06919"
$ws.Range("E12").Value = "  +0.79%  "
Set-TextValue "D13" "15.53"
$ws.Range("E13").Value = "  +5.46%  "
Set-TextValue "D14" "0.5992"
$ws.Range("E14").Value = "  -1.29%  "
Set-TextValue "D15" "4.489"
$ws.Range("E15").Value = "  +1.75%  "
Set-TextValue "D16" "77.15"
$ws.Range("E16").Value = "  -2.08%  "
$ws.Range("E17").Value = "  -0.73%  "
Set-TextValue "D18" "1.0000"
$ws.Range("E18").Value = "  -0.53%  "
Set-TextValue "D19" "25.773.32"
$ws.Range("E19").Value = "  -0.23%  "
Set-TextValue "D22" "1.971.17"
$ws.Range("E22").Value = "  -3.87%  "
Set-TextValue "D23" "4.073"
$ws.Range("E23").Value = "  +3.28%  "
Set-TextValue "D24" "8.256"
$ws.Range("E24").Value = "  +3.09%  "
Set-TextValue "D25" "5.193"
$ws.Range("E25").Value = "  -0.21%  "
Set-TextValue "D26" "136.81"
$ws.Range("E26").Value = "  +4.86%  "
Set-TextValue "D27" "1.446"
$ws.Range("E27").Value = "  +10.84%  "
Set-TextValue "D30" "102.28"
$ws.Range("E30").Value = "  +4.32%  "
Set-TextValue "D31" "0.08176"
$ws.Range("E31").Value = "  -1.35%  "
Set-TextValue "D32" "3.670"
$ws.Range("E32").Value = "  +2.03%  "
Set-TextValue "D33" "3.409"
$ws.Range("E33").Value = "  +8.27%  "
Set-TextValue "D34" "0.04495"
$ws.Range("E34").Value = "  +4.31%  "
$ws.Range("E35").Value = "  -0.53%  "
Set-TextValue "D36" "2.651"
$ws.Range("E36").Value = "  -3.29%  "
Set-TextValue "D37" "0.9982"
$ws.Range("E37").Value = "  -3.00%  "
Set-TextValue "D38" "0.5998"
$ws.Range("E38").Value = "  -3.78%  "
Set-TextValue "D39" "2.716"
$ws.Range("E39").Value = "  -5.35%  "
Set-TextValue "D40" "0.01554"
$ws.Range("E40").Value = "  +6.91%  "
Set-TextValue "D41" "1.921"
$ws.Range("E41").Value = "  -7.79%  "
$ws.Range("E42").Value = "  -0.50%  "
Set-TextValue "D43" "103.04"
$ws.Range("E43").Value = "  +2.62%  "
Set-TextValue "D44" "0.3785"
$ws.Range("E44").Value = "  -2.13%  "
Set-TextValue "D45" "0.7409"
$ws.Range("E45").Value = "  -5.25%  "
Set-TextValue "D46" "4.913"
$ws.Range("E46").Value = "  -4.60%  "
Set-TextValue "D47" "0.05480"
$ws.Range("E47").Value = "  +4.11%  "
Set-TextValue "D48" "0.1099"
$ws.Range("E48").Value = "  +7.46%  "
Set-TextValue "D49" "5.935"
$ws.Range("E49").Value = "  -4.08%  "
Set-TextValue "D50" "7.714"
$ws.Range("E50").Value = "  +2.90%  "
Set-TextValue "D51" "29.78"
$ws.Range("E51").Value = "  +2.73%  "

# Row 20/21 swap: Avalanche now ranks above ShibaInu, with updated price/volume
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D20" "11.64"
$ws.Range("E20").Value = "  +2.37%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D21" "0.000006799"
$ws.Range("E21").Value = "  +10.96%  "

# Row 28/29 swap: EthereumClassic now ranks above LidoDAOToken, with updated price/volume
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D28" "15.03"
$ws.Range("E28").Value = "  +4.15%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D29" "1.815"
$ws.Range("E29").Value = "  -2.98%  "
